# Insert a new "LP solver" row into the general sheet, right after the
# "NLP solver" row, shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("general")

# Insert a new row above current row 5 (Number of exp. conditions...)
$ws.Rows.Item(5).Insert()

$ws.Cells.Item(5, 1).Value = "LP solver (linprog or gurobi)"
$ws.Cells.Item(5, 2).Value = "gurobi"

# Copy the style from the label column of a neighboring row so the new
# label cell looks consistent with the rest of column A, then left-align
# the new label (distinguishing it slightly from the other label cells).
$ws.Cells.Item(4, 1).Copy()
$ws.Cells.Item(5, 1).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(5, 1).HorizontalAlignment = -4131  # xlLeft

$ws.Activate()
$ws.Range("A5:B5").Select()
